$d = $word.ActiveDocument

# The paragraph being edited is the only paragraph in the document body.
$p = $d.Paragraphs(1)
$r = $p.Range

$rPr = '<w:rPr><w:rFonts w:ascii="Helvetica" w:eastAsia="Times New Roman" w:hAnsi="Helvetica" w:cs="Arial"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr>'

$newParaInner = (
  '<w:r w:rsidRPr="001901AE">' + $rPr + '<w:t>This figure shows n</w:t></w:r>' +
  '<w:r w:rsidR="002F4DE4">' + $rPr + '<w:t xml:space="preserve">ational nutritional deficiencies by age and sex calculated using the Estimated Average Requirement (EAR) </w:t></w:r>' +
  '<w:proofErr w:type="spellStart"/>' +
  '<w:r>' + $rPr + '<w:t>cutpoint</w:t></w:r>' +
  '<w:proofErr w:type="spellEnd"/>' +
  '<w:r>' + $rPr + '<w:t xml:space="preserve"> method (National Academy</w:t></w:r>' +
  '<w:r w:rsidR="00792C57" w:rsidRPr="00792C57">' + $rPr + '<w:t>,</w:t></w:r>' +
  '<w:r w:rsidR="00792C57">' + $rPr + '<w:t xml:space="preserve"> 2000). Deeper orange shading indicates a higher proportion of individuals wit</w:t></w:r>' +
  '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>' +
  '<w:r>' + $rPr + '<w:t>hin an age-sex group that are deficient in each nutrient. The EAR is the average daily nutrient intake level estimated to meet the requirements of half of the healthy individuals in a group.</w:t></w:r>'
)

$pPr = '<w:pPr><w:rPr><w:rFonts w:ascii="Helvetica" w:eastAsia="Times New Roman" w:hAnsi="Helvetica" w:cs="Arial"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr></w:pPr>'

$newPara = '<w:p w14:paraId="6A95F2A6" w14:textId="36D241B3" w:rsidR="00E64A05" w:rsidRPr="00EC1F14" w:rsidRDefault="001901AE" w:rsidP="001901AE">' + $pPr + $newParaInner + '</w:p>'

$xml = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body>' + $newPara + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

# InsertXML replaces the whole paragraph (including its end-of-paragraph
# mark), which leaves one stray empty paragraph behind; remove the extra
# paragraph mark it introduces so the document keeps a single paragraph.
$r.InsertXML($xml)

$firstParaEnd = $d.Paragraphs(1).Range.End
$extraMark = $d.Range($firstParaEnd - 1, $firstParaEnd)
$extraMark.Delete()

Write-Output $d.Content.Text
